# Auto-generated Excel COM-interop script
# Applies numeric corrections to market-price / profit columns (H-N)
# across multiple rows on multiple crafting-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2094.4626
$ws.Range("I15").Value = 2094.4626
$ws.Range("K15").Value = 6283.3878
$ws.Range("M15").Value = -6114.3878

$ws.Range("H41").Value = 425.72726
$ws.Range("I41").Value = 268.4
$ws.Range("K41").Value = 268.4
$ws.Range("M41").Value = 171.6

$ws.Range("H64").Value = 4132.7334
$ws.Range("I64").Value = 2999.1
$ws.Range("J64").Value = 6400
$ws.Range("K64").Value = 2999.1
$ws.Range("L64").Value = 6400
$ws.Range("M64").Value = -2751.1
$ws.Range("N64").Value = -6896

$ws.Range("H67").Value = 4132.7334
$ws.Range("I67").Value = 2999.1
$ws.Range("J67").Value = 6400
$ws.Range("K67").Value = 2999.1
$ws.Range("L67").Value = 6400
$ws.Range("M67").Value = -2141.1
$ws.Range("N67").Value = -8116

$ws.Range("H69").Value = 7015
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 7015
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 21045
$ws.Range("N69").Value = -22793
$ws.Range("M69").ClearContents()

$ws.Range("H72").Value = 7015
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 7015
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 63135
$ws.Range("N72").Value = -71871
$ws.Range("M72").ClearContents()

$ws.Range("H86").Value = 10784.546
$ws.Range("J86").Value = 9998
$ws.Range("L86").Value = 9998
$ws.Range("N86").Value = -12244

$ws.Range("H89").Value = 10784.546
$ws.Range("J89").Value = 9998
$ws.Range("L89").Value = 49990
$ws.Range("N89").Value = -61222

$ws.Range("H100").Value = 2254.5881
$ws.Range("I100").Value = 1333
$ws.Range("K100").Value = 1333
$ws.Range("M100").Value = -792

$ws.Range("H131").Value = 30928.285
$ws.Range("I131").Value = 2874.5
$ws.Range("K131").Value = 8623.5
$ws.Range("M131").Value = -3583.5

$ws.Range("H138").Value = 4580.08
$ws.Range("J138").Value = 4933.925
$ws.Range("L138").Value = 14801.775
$ws.Range("N138").Value = -25081.775

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 191922
$ws.Range("I32").Value = 222298.05
$ws.Range("K32").Value = 222298.05
$ws.Range("M32").Value = -222011.05

$ws.Range("H97").Value = 37038940
$ws.Range("I97").Value = 1305.5
$ws.Range("K97").Value = 1305.5
$ws.Range("M97").Value = -809.5

$ws.Range("H132").Value = 1925248.4
$ws.Range("I132").Value = 2085352.4
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 6256057.199999999
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -6253527.199999999
$ws.Range("N132").Value = -17060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 6787
$ws.Range("I105").Value = 1555
$ws.Range("J105").Value = 8879.799999999999
$ws.Range("K105").Value = 1555
$ws.Range("L105").Value = 8879.799999999999
$ws.Range("M105").Value = 192
$ws.Range("N105").Value = -12373.8

$ws.Range("H107").Value = 3750
$ws.Range("I107").Value = 3500
$ws.Range("K107").Value = 3500
$ws.Range("M107").Value = -1580

$ws.Range("H135").Value = 76999.664
$ws.Range("J135").Value = 76999.664
$ws.Range("L135").Value = 76999.664
$ws.Range("N135").Value = -87139.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 65227.75
$ws.Range("I16").Value = 84470.336
$ws.Range("K16").Value = 84470.336
$ws.Range("M16").Value = -84183.336

$ws.Range("H99").Value = 35939.41
$ws.Range("I99").Value = 32227.277
$ws.Range("K99").Value = 32227.277
$ws.Range("M99").Value = -30729.277

$ws.Range("H105").Value = 5994.4546
$ws.Range("I105").Value = 6050.1665
$ws.Range("J105").Value = 5743.75
$ws.Range("K105").Value = 6050.1665
$ws.Range("L105").Value = 5743.75
$ws.Range("M105").Value = -4303.1665
$ws.Range("N105").Value = -9237.75

$ws.Range("H113").Value = 65227.75
$ws.Range("I113").Value = 84470.336
$ws.Range("K113").Value = 84470.336
$ws.Range("M113").Value = -82300.336

$ws.Range("H126").Value = 35939.41
$ws.Range("I126").Value = 32227.277
$ws.Range("K126").Value = 96681.83099999999
$ws.Range("M126").Value = -94211.83099999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 7886.9473
$ws.Range("I3").Value = 2685.3
$ws.Range("J3").Value = 13666.556
$ws.Range("K3").Value = 8055.900000000001
$ws.Range("L3").Value = 40999.66800000001
$ws.Range("M3").Value = -7943.900000000001
$ws.Range("N3").Value = -41223.66800000001

$ws.Range("H68").Value = 6176.25
$ws.Range("I68").Value = 1365.3334
$ws.Range("J68").Value = 7488.3184
$ws.Range("K68").Value = 4096.0002
$ws.Range("L68").Value = 22464.9552
$ws.Range("M68").Value = -3285.0002
$ws.Range("N68").Value = -24086.9552

$ws.Range("H71").Value = 6176.25
$ws.Range("I71").Value = 1365.3334
$ws.Range("J71").Value = 7488.3184
$ws.Range("K71").Value = 12288.0006
$ws.Range("L71").Value = 67394.8656
$ws.Range("M71").Value = -8232.000599999999
$ws.Range("N71").Value = -75506.8656

$ws.Range("H80").Value = 10000.6
$ws.Range("J80").Value = 10000.6
$ws.Range("L80").Value = 30001.8
$ws.Range("N80").Value = -31873.8

$ws.Range("H83").Value = 10000.6
$ws.Range("J83").Value = 10000.6
$ws.Range("L83").Value = 90005.40000000001
$ws.Range("N83").Value = -99365.40000000001

$ws.Range("H117").Value = 2278.6
$ws.Range("I117").Value = 1250
$ws.Range("J117").Value = 2964.3333
$ws.Range("K117").Value = 3750
$ws.Range("L117").Value = 8892.999899999999
$ws.Range("M117").Value = -308
$ws.Range("N117").Value = -15776.9999

$ws.Range("H122").Value = 897749.4
$ws.Range("J122").Value = 1973.0769
$ws.Range("L122").Value = 17757.6921
$ws.Range("N122").Value = -22657.6921

$ws.Range("H131").Value = 4972.276
$ws.Range("J131").Value = 5747.913
$ws.Range("L131").Value = 17243.739
$ws.Range("N131").Value = -27323.739

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1504.6875
$ws.Range("I102").Value = 813.4167
$ws.Range("K102").Value = 813.4167
$ws.Range("M102").Value = 808.5833

$ws.Range("H113").Value = 4917.143
$ws.Range("I113").Value = 3886
$ws.Range("K113").Value = 3886
$ws.Range("M113").Value = -1716

$ws.Range("H132").Value = 22293.818
$ws.Range("I132").Value = 20058.334
$ws.Range("J132").Value = 24976.4
$ws.Range("K132").Value = 60175.00199999999
$ws.Range("L132").Value = 74929.20000000001
$ws.Range("M132").Value = -57645.00199999999
$ws.Range("N132").Value = -79989.20000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5381.1
$ws.Range("I46").Value = 1277.3334
$ws.Range("K46").Value = 1277.3334
$ws.Range("M46").Value = -1089.3334

$ws.Range("H132").Value = 5575735
$ws.Range("I132").Value = 11145724
$ws.Range("J132").Value = 5745.8335
$ws.Range("K132").Value = 33437172
$ws.Range("L132").Value = 17237.5005
$ws.Range("M132").Value = -33434642
$ws.Range("N132").Value = -22297.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 31782.428
$ws.Range("J62").Value = 31919.54
$ws.Range("L62").Value = 31919.54
$ws.Range("N62").Value = -33167.54

$ws.Range("H65").Value = 31782.428
$ws.Range("J65").Value = 31919.54
$ws.Range("L65").Value = 159597.7
$ws.Range("N65").Value = -165837.7

$ws.Range("H96").Value = 1755.4445
$ws.Range("J96").Value = 1899.75
$ws.Range("L96").Value = 1899.75
$ws.Range("N96").Value = -4645.75

$ws.Range("H100").Value = 1999.25
$ws.Range("I100").Value = 1999.25
$ws.Range("K100").Value = 3998.5
$ws.Range("M100").Value = -3457.5

$ws.Range("H132").Value = 5378859
$ws.Range("I132").Value = 6175331.5
$ws.Range("J132").Value = 2668.75
$ws.Range("K132").Value = 18525994.5
$ws.Range("L132").Value = 8006.25
$ws.Range("M132").Value = -18523464.5
$ws.Range("N132").Value = -13066.25
